$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 10 ("Greater than USD100 CR" rule), shifting subsequent rows up
$ws.Rows.Item(10).Delete()

# Update selection to match the post-edit state (A10:XFD10)
$ws.Rows.Item(10).Select()
